# Update cryptocurrency price/volume data (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '27.255.19'
$ws.Cells.Item(2, 5).Value = '  +1.15%  '

$ws.Cells.Item(3, 4).Value = '1.569.92'
$ws.Cells.Item(3, 5).Value = '  +0.65%  '

$ws.Cells.Item(4, 5).Value = '  +0.16%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '210.91'
$ws.Cells.Item(5, 5).Value = '  +1.82%  '

$ws.Cells.Item(6, 5).Value = '  +0.56%  '

$ws.Cells.Item(7, 5).Value = '  +0.18%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '22.07'
$ws.Cells.Item(8, 5).Value = '  +0.01%  '

$ws.Cells.Item(9, 5).Value = '  +0.35%  '

$ws.Cells.Item(10, 5).Value = '  -0.07%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0869'
$ws.Cells.Item(11, 5).Value = '  +1.40%  '

$ws.Cells.Item(12, 4).Value = '1.794.11'
$ws.Cells.Item(12, 5).Value = '  +0.67%  '

$ws.Cells.Item(13, 4).Value = '1.559.54'
$ws.Cells.Item(13, 5).Value = '  -0.02%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '3.78'
$ws.Cells.Item(14, 5).Value = '  +0.70%  '

$ws.Cells.Item(15, 5).Value = '  +0.14%  '

$ws.Cells.Item(16, 4).Value = '27.194.44'
$ws.Cells.Item(16, 5).Value = '  +0.90%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '62.21'
$ws.Cells.Item(17, 5).Value = '  +0.23%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '7.56'
$ws.Cells.Item(18, 5).Value = '  +2.90%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '216.22'
$ws.Cells.Item(19, 5).Value = '  -0.25%  '

$ws.Cells.Item(20, 5).Value = '  -0.39%  '

$ws.Cells.Item(21, 5).Value = '  +0.21%  '

$ws.Cells.Item(22, 5).Value = '  +1.08%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '9.22'
$ws.Cells.Item(23, 5).Value = '  +0.15%  '

$ws.Cells.Item(24, 5).Value = '  +0.24%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '153.82'
$ws.Cells.Item(25, 5).Value = '  +0.74%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '6.63'
$ws.Cells.Item(26, 5).Value = '  +0.49%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '15.07'
$ws.Cells.Item(27, 5).Value = '  +0.07%  '

$ws.Cells.Item(28, 5).Value = '  +2.25%  '

$ws.Cells.Item(29, 5).Value = '  +0.24%  '

$ws.Cells.Item(30, 5).Value = '  +2.52%  '

$ws.Cells.Item(31, 5).Value = '  +0.43%  '

$ws.Cells.Item(32, 5).Value = '  +0.13%  '

$ws.Cells.Item(33, 4).Value = '1.449.26'
$ws.Cells.Item(33, 5).Value = '  +2.23%  '

$ws.Cells.Item(34, 5).Value = '  +1.40%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.13'
$ws.Cells.Item(35, 5).Value = '  +6.91%  '

$ws.Cells.Item(36, 5).Value = '  +0.26%  '

$ws.Cells.Item(37, 5).Value = '  +0.37%  '

$ws.Cells.Item(38, 5).Value = '  +0.99%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.533'
$ws.Cells.Item(39, 5).Value = '  +0.22%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '5.88'
$ws.Cells.Item(40, 5).Value = '  +2.73%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.810'
$ws.Cells.Item(41, 5).Value = '  +0.29%  '

$ws.Cells.Item(42, 5).Value = '  +0.23%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '2.34'
$ws.Cells.Item(43, 5).Value = '  +0.75%  '

$ws.Cells.Item(44, 5).Value = '  -0.66%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '64.40'
$ws.Cells.Item(45, 5).Value = '  -0.50%  '

$ws.Cells.Item(46, 5).Value = '  -1.21%  '

$ws.Cells.Item(47, 4).Value = '1.706.01'
$ws.Cells.Item(47, 5).Value = '  +0.56%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '86.08'
$ws.Cells.Item(48, 5).Value = '  -1.54%  '

$ws.Cells.Item(49, 5).Value = '  +0.76%  '

$ws.Cells.Item(50, 5).Value = '  +1.60%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.0954'
$ws.Cells.Item(51, 5).Value = '  -0.15%  '
